$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OverallRebateEfficiency")
$ws.Range("A1").Value = "Week"
